$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.169.85"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.856.82"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.63"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4688"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2889"
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06558"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.83"
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07976"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.44"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "1.855.23"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.104"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6772"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.44"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").Value = "30.145.38"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("E18").Value = "  +7.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007654"
$ws.Range("E19").Value = "  +4.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9991"
$ws.Range("D21").Value = "2.097.13"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9982"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.207"
$ws.Range("E23").Value = "  -4.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.151"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.85"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.177"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.93"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.939"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.380"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09881"
$ws.Range("E30").Value = "  +2.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.467"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.311"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.024"
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04706"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.120"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6986"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01870"
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.605"
$ws.Range("E39").Value = "  +2.74%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.38"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.932"
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8395"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9982"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.46"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4143"
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.161"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "937.29"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.040"
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.94"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05651"
